$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: SCD0179 -> SCD0011
$ws.Name = "SCD0011"

# Update the TC_ID cell content (B2/B3): DGS-194 -> SCD0011-010
$ws.Range("B2").Value = "SCD0011-010"
$ws.Range("B3").Value = "SCD0011-010"

# Column B widened to fit the new (longer) TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.592447916666666

# Move the active selection from E2 to B4 (and drop the old frozen/top-left anchor at D1)
$ws.Range("B4").Select() | Out-Null

# Best-effort: reposition/resize the workbook window to match the saved view
$win = $excel.ActiveWindow
$win.Left = -120
$win.Top = -120
$win.Width = 20730
$win.Height = 11760
